$d = $word.ActiveDocument

# --- Locate the "Huis en Inrichting" heading paragraph (top-level category) ---
$heading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Huis en Inrichting") {
        $heading = $p
        break
    }
}

# --- Insert a new subcategory paragraph "Huizen" right after the heading ---
$newPara = $heading.Range.InsertParagraphAfter()
$huizenPara = $heading.Next()
$huizenPara.Range.ListFormat.ListLevelNumber = 2
$huizenPara.Range.Text = "Huizen"
$huizenPara.Range.Font.Bold = $true

# --- Move the "_GoBack" bookmark to the start of the new "Huizen" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$huizenPara = $heading.Next()
$bmRange = $d.Range($huizenPara.Range.Start, $huizenPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Remove the old top-level "Huizen" paragraph (it sat between the old
#     "Overige" sub-item of "Huis en Inrichting" and "Kinderen en Baby's") ---
$oldHuizen = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Huizen" -and
        $p.Range.ListFormat.ListLevelNumber -eq 1) {
        $oldHuizen = $p
        break
    }
}
$oldHuizen.Range.Delete()

Write-Output "done"
